# Commit: "Analises por regiao, regressao por regiao"
# This removes Portuguese accents/diacritics from the municipality names
# (stored as shared strings referenced from column A of every sheet) so
# the data can be used for region-based analyses without encoding issues.

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item(1)  # "info"
$ws4Serie = $wb.Worksheets.Item(2)  # "4 serie5 ano"
$ws8Serie = $wb.Worksheets.Item(3)  # "8 Serie 9 ano"
$ws3EM = $wb.Worksheets.Item(4)  # "3 serie EM"

# --- Strip accents from municipality names on "4 serie5 ano" (94 cells) ---
$ws4SerieUpdates = @(
  @("A1", "Municipio"),
  @("A2", "Caldas Brandao"),
  @("A5", "agua Branca"),
  @("A6", "Esperanca"),
  @("A9", "Olho dagua"),
  @("A10", "Sao Mamede"),
  @("A15", "Bonito de Santa Fe"),
  @("A16", "Dona Ines"),
  @("A18", "Marizopolis"),
  @("A19", "Santa Ines"),
  @("A20", "Caturite"),
  @("A23", "Tenorio"),
  @("A24", "Zabele"),
  @("A26", "Camalau"),
  @("A28", "Catole do Rocha"),
  @("A30", "Sao Bentinho"),
  @("A31", "Baia da Traicao"),
  @("A39", "Sao Joao do Cariri"),
  @("A41", "Caapora"),
  @("A43", "Joao Pessoa"),
  @("A44", "Puxinana"),
  @("A46", "Sao Sebastiao de Lagoa de Roca"),
  @("A57", "Sao Francisco"),
  @("A58", "Sao Miguel de Taipu"),
  @("A59", "Sume"),
  @("A62", "Conceicao"),
  @("A63", "Cruz do Espirito Santo"),
  @("A76", "Boqueirao"),
  @("A82", "Jerico"),
  @("A86", "Picui"),
  @("A88", "Sao Jose do Brejo do Cruz"),
  @("A89", "Sao Sebastiao do Umbuzeiro"),
  @("A90", "Uirauna"),
  @("A92", "Barauna"),
  @("A93", "Barra de Sao Miguel"),
  @("A95", "Belem"),
  @("A97", "Caicara"),
  @("A102", "Riachao"),
  @("A104", "Sao Jose do Sabugi"),
  @("A109", "Damiao"),
  @("A110", "Junco do Serido"),
  @("A112", "Piloezinhos"),
  @("A117", "Sao Jose do Bonfim"),
  @("A119", "Solanea"),
  @("A122", "Juarez Tavora"),
  @("A123", "Matureia"),
  @("A126", "Remigio"),
  @("A131", "Cachoeira dos indios"),
  @("A133", "Cuite de Mamanguape"),
  @("A134", "Inga"),
  @("A135", "Jacarau"),
  @("A136", "Mae dagua"),
  @("A140", "Riachao do Bacamarte"),
  @("A141", "Sao Jose de Piranhas"),
  @("A144", "Sossego"),
  @("A148", "Assuncao"),
  @("A151", "Cuite"),
  @("A153", "Manaira"),
  @("A157", "Pianco"),
  @("A158", "Poco de Jose de Moura"),
  @("A159", "Sao Jose dos Ramos"),
  @("A160", "Sape"),
  @("A161", "Aracagi"),
  @("A163", "Gurinhem"),
  @("A165", "Riachao do Poco"),
  @("A168", "Belem do Brejo do Cruz"),
  @("A170", "Sao Jose da Lagoa Tapada"),
  @("A174", "Marcacao"),
  @("A177", "Sao Bento"),
  @("A178", "Piloes"),
  @("A180", "Algodao de Jandaira"),
  @("A182", "Areia de Baraunas"),
  @("A185", "Caraubas"),
  @("A190", "Gurjao"),
  @("A202", "Pedro Regis"),
  @("A203", "Poco Dantas"),
  @("A204", "Quixaba"),
  @("A205", "Riacho de Santo Antonio"),
  @("A207", "Salgado de Sao Felix"),
  @("A208", "Santa Cecilia"),
  @("A210", "Santo Andre"),
  @("A211", "Sao Domingos"),
  @("A212", "Sao Domingos do Cariri"),
  @("A213", "Sao Joao do Rio do Peixe"),
  @("A214", "Sao Joao do Tigre"),
  @("A215", "Sao Jose de Caiana"),
  @("A216", "Sao Jose de Espinharas"),
  @("A217", "Sao Jose de Princesa"),
  @("A218", "Sao Jose dos Cordeiros"),
  @("A219", "Sao Vicente do Serido"),
  @("A220", "Sertaozinho"),
  @("A221", "Taperoa"),
  @("A222", "Varzea"),
  @("A223", "Vieiropolis")
)
foreach ($pair in $ws4SerieUpdates) {
  $ws4Serie.Range($pair[0]).Value = $pair[1]
}

# --- Strip accents from municipality names on "8 Serie 9 ano" (92 cells) ---
$ws8SerieUpdates = @(
  @("A1", "Municipio"),
  @("A2", "agua Branca"),
  @("A8", "Algodao de Jandaira"),
  @("A12", "Aracagi"),
  @("A16", "Areia de Baraunas"),
  @("A19", "Assuncao"),
  @("A20", "Baia da Traicao"),
  @("A22", "Barauna"),
  @("A25", "Barra de Sao Miguel"),
  @("A27", "Belem"),
  @("A28", "Belem do Brejo do Cruz"),
  @("A34", "Bonito de Santa Fe"),
  @("A35", "Boqueirao"),
  @("A39", "Caapora"),
  @("A42", "Cachoeira dos indios"),
  @("A46", "Caicara"),
  @("A49", "Caldas Brandao"),
  @("A50", "Camalau"),
  @("A53", "Caraubas"),
  @("A57", "Catole do Rocha"),
  @("A58", "Caturite"),
  @("A59", "Conceicao"),
  @("A64", "Cruz do Espirito Santo"),
  @("A66", "Cuite"),
  @("A67", "Cuite de Mamanguape"),
  @("A71", "Damiao"),
  @("A74", "Dona Ines"),
  @("A77", "Esperanca"),
  @("A82", "Gurinhem"),
  @("A83", "Gurjao"),
  @("A87", "Inga"),
  @("A92", "Jacarau"),
  @("A93", "Jerico"),
  @("A94", "Joao Pessoa"),
  @("A96", "Juarez Tavora"),
  @("A98", "Junco do Serido"),
  @("A108", "Mae dagua"),
  @("A111", "Manaira"),
  @("A112", "Marcacao"),
  @("A114", "Marizopolis"),
  @("A119", "Matureia"),
  @("A130", "Olho dagua"),
  @("A139", "Pedro Regis"),
  @("A140", "Pianco"),
  @("A141", "Picui"),
  @("A143", "Piloes"),
  @("A144", "Piloezinhos"),
  @("A148", "Poco Dantas"),
  @("A149", "Poco de Jose de Moura"),
  @("A153", "Puxinana"),
  @("A155", "Quixaba"),
  @("A156", "Remigio"),
  @("A157", "Riachao"),
  @("A158", "Riachao do Bacamarte"),
  @("A159", "Riachao do Poco"),
  @("A160", "Riacho de Santo Antonio"),
  @("A164", "Salgado de Sao Felix"),
  @("A165", "Santa Cecilia"),
  @("A168", "Santa Ines"),
  @("A174", "Santo Andre"),
  @("A175", "Sao Bentinho"),
  @("A176", "Sao Bento"),
  @("A177", "Sao Domingos"),
  @("A178", "Sao Domingos do Cariri"),
  @("A179", "Sao Francisco"),
  @("A180", "Sao Joao do Cariri"),
  @("A181", "Sao Joao do Rio do Peixe"),
  @("A182", "Sao Joao do Tigre"),
  @("A183", "Sao Jose da Lagoa Tapada"),
  @("A184", "Sao Jose de Caiana"),
  @("A185", "Sao Jose de Espinharas"),
  @("A186", "Sao Jose de Piranhas"),
  @("A187", "Sao Jose do Bonfim"),
  @("A188", "Sao Jose do Brejo do Cruz"),
  @("A189", "Sao Jose do Sabugi"),
  @("A190", "Sao Jose dos Ramos"),
  @("A191", "Sao Mamede"),
  @("A192", "Sao Miguel de Taipu"),
  @("A193", "Sao Sebastiao de Lagoa de Roca"),
  @("A194", "Sao Sebastiao do Umbuzeiro"),
  @("A195", "Sao Vicente do Serido"),
  @("A196", "Sape"),
  @("A202", "Sertaozinho"),
  @("A204", "Solanea"),
  @("A206", "Sossego"),
  @("A208", "Sume"),
  @("A210", "Taperoa"),
  @("A213", "Tenorio"),
  @("A215", "Uirauna"),
  @("A217", "Varzea"),
  @("A218", "Vieiropolis"),
  @("A219", "Zabele")
)
foreach ($pair in $ws8SerieUpdates) {
  $ws8Serie.Range($pair[0]).Value = $pair[1]
}

# --- Strip accents from municipality names on "3 serie EM" (2 cells) ---
$ws3EMUpdates = @(
  @("A2", "Municipio"),
  @("A4", "Sao Domingos do Cariri")
)
foreach ($pair in $ws3EMUpdates) {
  $ws3EM.Range($pair[0]).Value = $pair[1]
}

# --- Restore sheet selections / active sheet ---
# Originally "8 Serie 9 ano" (index 3) was the active tab (activeTab=2, 0-based);
# now "4 serie5 ano" (index 2) is the active tab (activeTab=1, 0-based).
# Activation order matters: activate each sheet, then finish on the one that
# should end up as the active tab.

[void]$wsInfo.Activate()
[void]$wsInfo.Range("A1").Select()

[void]$ws3EM.Activate()
[void]$ws3EM.Range("A1").Select()

[void]$ws8Serie.Activate()
[void]$ws8Serie.Range("A108").Select()

[void]$ws4Serie.Activate()
[void]$ws4Serie.Range("A26").Select()

